$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: row 4 (6e85fe82-... file) handoff/handback datetimes are regenerated
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-17 16:10:50"
$wsZhCn.Range("G4").Value = "2016-01-17 16:11:37"

# "de-de" sheet: row 4 (6e85fe82-... file) handoff/handback datetimes are regenerated
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-17 16:11:01"
$wsDeDe.Range("G4").Value = "2016-01-17 16:11:57"
